$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 316, shifting existing rows 316:384 down to 317:385
$ws.Rows("316:316").Insert()

# Populate the newly inserted row 316 with the new record's data
$ws.Range("A316").Value = 3
$ws.Range("B316").Value = "Femacal de La Calera"
$ws.Range("C316").Value = "Coquimbo"
$ws.Range("D316").Value = 45218
$ws.Range("E316").Value = 5
$ws.Range("F316").Value = "Fruta"
$ws.Range("G316").Value = 100101
$ws.Range("H316").Value = "Berries"
$ws.Range("I316").Value = 100101001
$ws.Range("J316").Value = "Arándano (blue)"
$ws.Range("K316").Value = "Sin especificar"
$ws.Range("L316").Value = "Primera"
$ws.Range("M316").Value = 56
$ws.Range("N316").Value = 13000
$ws.Range("O316").Value = 13000
$ws.Range("P316").Value = 13000
$ws.Range("Q316").Value = "$/bandeja 2 kilos"
$ws.Range("R316").Value = "Provincia de Quillota"
$ws.Range("S316").Value = 6500
$ws.Range("T316").Value = 2
